$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 0.05871566666666667
$ws.Range("N2").Value = 0.176147
$ws.Range("O2").Value = 0.008355169877275808
$ws.Range("P2").Value = 0.008355169877275808
$ws.Range("Q2").Value = 0.02326369514622222
$ws.Range("R2").Value = 0.209373256316
$ws.Range("S2").Value = 0.008355169877275808
$ws.Range("T2").Value = 0.008355169877275808

# Row 3 updates
$ws.Range("O3").Value = 0.1868088427899751
$ws.Range("P3").Value = 0.1868088427899751
$ws.Range("S3").Value = 0.1868088427899751
$ws.Range("T3").Value = 0.1868088427899751

# Row 4 updates
$ws.Range("O4").Value = 0.8048359873327491
$ws.Range("P4").Value = 0.8048359873327491
$ws.Range("S4").Value = 0.8048359873327491
$ws.Range("T4").Value = 0.8048359873327491
